# Split Project1 into Project1a and Project1b
#
# This script edits slide 2 ("System Overview - v1.0 (Project 2)") of the
# presentation: it removes the enclosing background rounded-rectangle and
# the second MongoDB/Database "stack" (and its trailing dot-leader textbox),
# moves the remaining MongoDB/Database box down so it lines up on its own,
# and re-routes the two connector arrows that used to fan out to both
# stacks so they both point at the single remaining box.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

function Get-ShapeById {
    param($slide, $id)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

# 1. Remove the big background rounded rectangle that used to sit behind
#    both MongoDB/Database boxes (id 30, "Rectangle: Rounded Corners 29").
(Get-ShapeById $s 30).Delete()

# 2. Move the remaining MongoDB/Database box (id 6) down into the spot
#    previously centered between the two stacks it now replaces.
$mongo = Get-ShapeById $s 6
$mongo.Left = 776.7341732283464
$mongo.Top = 243.04378512755906

# 3. Remove the second MongoDB/Database box (id 31, "Rectangle: Rounded
#    Corners 30") - it is being merged into the single box handled above.
(Get-ShapeById $s 31).Delete()

# 4. Remove the orphaned dot-leader textbox that sat to the right of the
#    box removed above (id 33, "TextBox 32").
(Get-ShapeById $s 33).Delete()

# 5. Re-point the first connector (id 34) so it travels from the REST box
#    down to the (moved) MongoDB/Database box - it no longer needs to be
#    flipped vertically since it now travels downward on-screen.
$conn34 = Get-ShapeById $s 34
$conn34.Left = 654.640167280315
$conn34.Top = 203.81464566929134
$conn34.Width = 122.0940157480315
$conn34.Height = 69.7244094488189
$conn34.VerticalFlip = $false

# 6. Re-point the second connector (id 36) so it also travels up to the
#    single remaining MongoDB/Database box, and glue its end to it.
$conn36 = Get-ShapeById $s 36
$conn36.Left = 653.9993700787402
$conn36.Top = 273.53905511811024
$conn36.Width = 122.7348061496063
$conn36.Height = 69.72433090866141
$conn36.ConnectorFormat.EndConnect($mongo, 1)
